# Auto-generated edit script applying the Sargatanas_Profits.xlsx numeric diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 883.6667
$ws.Cells.Item(4, 9).Value = 883.6667
$ws.Cells.Item(4, 11).Value = 883.6667
$ws.Cells.Item(4, 13).Value = -769.6667
$ws.Cells.Item(6, 8).Value = 34482760
$ws.Cells.Item(6, 9).Value = 34482760
$ws.Cells.Item(6, 11).Value = 103448280
$ws.Cells.Item(6, 13).Value = -103448168
$ws.Cells.Item(43, 8).Value = 1399996.4
$ws.Cells.Item(43, 9).Value = 49994.75
$ws.Cells.Item(43, 11).Value = 49994.75
$ws.Cells.Item(43, 13).Value = -49925.75
$ws.Cells.Item(62, 8).Value = 53192.5
$ws.Cells.Item(62, 9).Value = 10000
$ws.Cells.Item(62, 11).Value = 10000
$ws.Cells.Item(62, 13).Value = -9376
$ws.Cells.Item(64, 8).Value = 7981.6665
$ws.Cells.Item(64, 10).Value = 7981.6665
$ws.Cells.Item(64, 12).Value = 7981.6665
$ws.Cells.Item(64, 14).Value = -8477.666499999999
$ws.Cells.Item(65, 8).Value = 53192.5
$ws.Cells.Item(65, 9).Value = 10000
$ws.Cells.Item(65, 11).Value = 50000
$ws.Cells.Item(65, 13).Value = -46880
$ws.Cells.Item(67, 8).Value = 7981.6665
$ws.Cells.Item(67, 10).Value = 7981.6665
$ws.Cells.Item(67, 12).Value = 7981.6665
$ws.Cells.Item(67, 14).Value = -9697.666499999999
$ws.Cells.Item(69, 8).Value = 6000
$ws.Cells.Item(69, 10).Value = 6000
$ws.Cells.Item(69, 12).Value = 18000
$ws.Cells.Item(69, 14).Value = -19748
$ws.Cells.Item(72, 8).Value = 6000
$ws.Cells.Item(72, 10).Value = 6000
$ws.Cells.Item(72, 12).Value = 54000
$ws.Cells.Item(72, 14).Value = -62736
$ws.Cells.Item(98, 8).Value = 3440.7778
$ws.Cells.Item(98, 9).Value = 3489.1667
$ws.Cells.Item(98, 10).Value = 3198.8333
$ws.Cells.Item(98, 11).Value = 3489.1667
$ws.Cells.Item(98, 12).Value = 3198.8333
$ws.Cells.Item(98, 13).Value = -1991.1667
$ws.Cells.Item(98, 14).Value = -6194.8333
$ws.Cells.Item(105, 8).Value = 54797
$ws.Cells.Item(105, 10).Value = 54797
$ws.Cells.Item(105, 12).Value = 54797
$ws.Cells.Item(105, 14).Value = -61785
$ws.Cells.Item(106, 8).Value = 1900.9445
$ws.Cells.Item(106, 9).Value = 1900.9445
$ws.Cells.Item(106, 11).Value = 1900.9445
$ws.Cells.Item(106, 13).Value = -1269.9445
$ws.Cells.Item(107, 8).Value = 21593294
$ws.Cells.Item(107, 9).Value = 8335683.5
$ws.Cells.Item(107, 10).Value = 50002456
$ws.Cells.Item(107, 11).Value = 8335683.5
$ws.Cells.Item(107, 12).Value = 50002456
$ws.Cells.Item(107, 13).Value = -8333763.5
$ws.Cells.Item(107, 14).Value = -50006296
$ws.Cells.Item(113, 8).Value = 31258250
$ws.Cells.Item(113, 9).Value = 5562.375
$ws.Cells.Item(113, 10).Value = 62510940
$ws.Cells.Item(113, 11).Value = 5562.375
$ws.Cells.Item(113, 12).Value = 62510940
$ws.Cells.Item(113, 13).Value = -2308.375
$ws.Cells.Item(113, 14).Value = -62517448
$ws.Cells.Item(118, 8).Value = 1874.1428
$ws.Cells.Item(118, 9).Value = 485
$ws.Cells.Item(118, 11).Value = 1455
$ws.Cells.Item(118, 13).Value = 202
$ws.Cells.Item(122, 8).Value = 3440.7778
$ws.Cells.Item(122, 9).Value = 3489.1667
$ws.Cells.Item(122, 10).Value = 3198.8333
$ws.Cells.Item(122, 11).Value = 10467.5001
$ws.Cells.Item(122, 12).Value = 9596.499899999999
$ws.Cells.Item(122, 13).Value = -8017.500100000001
$ws.Cells.Item(122, 14).Value = -14496.4999
$ws.Cells.Item(129, 8).Value = 1321.5454
$ws.Cells.Item(129, 10).Value = 2428.8
$ws.Cells.Item(129, 12).Value = 7286.400000000001
$ws.Cells.Item(129, 14).Value = -17286.4
$ws.Cells.Item(132, 8).Value = 1154.1111
$ws.Cells.Item(132, 9).Value = 1106.279
$ws.Cells.Item(132, 10).Value = 2182.5
$ws.Cells.Item(132, 11).Value = 3318.837
$ws.Cells.Item(132, 12).Value = 6547.5
$ws.Cells.Item(132, 13).Value = -788.837
$ws.Cells.Item(132, 14).Value = -11607.5
$ws.Cells.Item(135, 8).Value = 400936.88
$ws.Cells.Item(135, 9).Value = 417600.97
$ws.Cells.Item(135, 10).Value = 999
$ws.Cells.Item(135, 11).Value = 3758408.73
$ws.Cells.Item(135, 12).Value = 8991
$ws.Cells.Item(135, 13).Value = -3755873.73
$ws.Cells.Item(135, 14).Value = -14061
$ws.Cells.Item(137, 8).Value = 2058.45
$ws.Cells.Item(137, 9).Value = 1756.2727
$ws.Cells.Item(137, 10).Value = 2427.7778
$ws.Cells.Item(137, 11).Value = 5268.8181
$ws.Cells.Item(137, 12).Value = 7283.3334
$ws.Cells.Item(137, 13).Value = -2718.8181
$ws.Cells.Item(137, 14).Value = -12383.3334
$ws.Cells.Item(138, 8).Value = 3849487.2
$ws.Cells.Item(138, 10).Value = 8338224
$ws.Cells.Item(138, 12).Value = 25014672
$ws.Cells.Item(138, 14).Value = -25024952
$ws.Cells.Item(141, 8).Value = 2479.4285
$ws.Cells.Item(141, 9).Value = 2478.4
$ws.Cells.Item(141, 11).Value = 7435.200000000001
$ws.Cells.Item(141, 13).Value = -2255.200000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 10).Value = 10000
$ws.Cells.Item(10, 12).Value = 10000
$ws.Cells.Item(10, 14).Value = -10340
$ws.Cells.Item(32, 8).Value = 1646898.6
$ws.Cells.Item(32, 9).Value = 1787922.1
$ws.Cells.Item(32, 10).Value = 1623.8334
$ws.Cells.Item(32, 11).Value = 1787922.1
$ws.Cells.Item(32, 12).Value = 1623.8334
$ws.Cells.Item(32, 13).Value = -1787635.1
$ws.Cells.Item(32, 14).Value = -2197.8334
$ws.Cells.Item(38, 8).Value = 32205.8
$ws.Cells.Item(38, 9).Value = 32205.8
$ws.Cells.Item(38, 11).Value = 32205.8
$ws.Cells.Item(38, 13).Value = -31738.8
$ws.Cells.Item(45, 8).Value = 5337
$ws.Cells.Item(45, 9).Value = 2276
$ws.Cells.Item(45, 11).Value = 2276
$ws.Cells.Item(45, 13).Value = -1899
$ws.Cells.Item(61, 8).Value = 5970.8
$ws.Cells.Item(61, 9).Value = 2746.3215
$ws.Cells.Item(61, 10).Value = 10074.682
$ws.Cells.Item(61, 11).Value = 2746.3215
$ws.Cells.Item(61, 12).Value = 10074.682
$ws.Cells.Item(61, 13).Value = -2534.3215
$ws.Cells.Item(61, 14).Value = -10498.682
$ws.Cells.Item(62, 8).Value = 42999.5
$ws.Cells.Item(62, 10).Value = 42999.5
$ws.Cells.Item(62, 12).Value = 42999.5
$ws.Cells.Item(62, 14).Value = -44247.5
$ws.Cells.Item(65, 8).Value = 42999.5
$ws.Cells.Item(65, 10).Value = 42999.5
$ws.Cells.Item(65, 12).Value = 128998.5
$ws.Cells.Item(65, 14).Value = -135238.5
$ws.Cells.Item(74, 8).Value = 46783.832
$ws.Cells.Item(74, 9).Value = 66188.44
$ws.Cells.Item(74, 10).Value = 7974.625
$ws.Cells.Item(74, 11).Value = 66188.44
$ws.Cells.Item(74, 12).Value = 7974.625
$ws.Cells.Item(74, 13).Value = -65314.44
$ws.Cells.Item(74, 14).Value = -9722.625
$ws.Cells.Item(77, 8).Value = 46783.832
$ws.Cells.Item(77, 9).Value = 66188.44
$ws.Cells.Item(77, 10).Value = 7974.625
$ws.Cells.Item(77, 11).Value = 330942.2
$ws.Cells.Item(77, 12).Value = 39873.125
$ws.Cells.Item(77, 13).Value = -326574.2
$ws.Cells.Item(77, 14).Value = -48609.125
$ws.Cells.Item(102, 8).Value = 1294.6923
$ws.Cells.Item(102, 9).Value = 1327.5
$ws.Cells.Item(102, 11).Value = 1327.5
$ws.Cells.Item(102, 13).Value = 294.5
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 14).Value = $null
$ws.Cells.Item(110, 8).Value = 25642014
$ws.Cells.Item(110, 9).Value = 798.875
$ws.Cells.Item(110, 10).Value = 66667956
$ws.Cells.Item(110, 11).Value = 798.875
$ws.Cells.Item(110, 12).Value = 66667956
$ws.Cells.Item(110, 13).Value = 1246.125
$ws.Cells.Item(110, 14).Value = -66672046
$ws.Cells.Item(122, 8).Value = 1712
$ws.Cells.Item(122, 9).Value = 1302.3334
$ws.Cells.Item(122, 11).Value = 3907.0002
$ws.Cells.Item(122, 13).Value = -1457.0002
$ws.Cells.Item(132, 8).Value = 1004980.75
$ws.Cells.Item(132, 9).Value = 1474223.2
$ws.Cells.Item(132, 10).Value = 7840.5625
$ws.Cells.Item(132, 11).Value = 4422669.6
$ws.Cells.Item(132, 12).Value = 23521.6875
$ws.Cells.Item(132, 13).Value = -4420139.6
$ws.Cells.Item(132, 14).Value = -28581.6875
$ws.Cells.Item(134, 8).Value = 53749.5
$ws.Cells.Item(134, 10).Value = 53749.5
$ws.Cells.Item(134, 12).Value = 53749.5
$ws.Cells.Item(134, 14).Value = -63889.5
$ws.Cells.Item(136, 8).Value = 5970.8
$ws.Cells.Item(136, 9).Value = 2746.3215
$ws.Cells.Item(136, 10).Value = 10074.682
$ws.Cells.Item(136, 11).Value = 8238.9645
$ws.Cells.Item(136, 12).Value = 30224.046
$ws.Cells.Item(136, 13).Value = -5688.9645
$ws.Cells.Item(136, 14).Value = -35324.046
$ws.Cells.Item(140, 8).Value = 74888.5
$ws.Cells.Item(140, 10).Value = 74888.5
$ws.Cells.Item(140, 12).Value = 74888.5
$ws.Cells.Item(140, 14).Value = -85248.5
$ws.Cells.Item(141, 8).Value = 26476.334
$ws.Cells.Item(141, 10).Value = 26476.334
$ws.Cells.Item(141, 12).Value = 26476.334
$ws.Cells.Item(141, 14).Value = -36836.334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 15153761
$ws.Cells.Item(20, 9).Value = 23810994
$ws.Cells.Item(20, 11).Value = 23810994
$ws.Cells.Item(20, 13).Value = -23810747
$ws.Cells.Item(62, 8).Value = 43000
$ws.Cells.Item(62, 10).Value = 43000
$ws.Cells.Item(62, 12).Value = 43000
$ws.Cells.Item(62, 14).Value = -44372
$ws.Cells.Item(65, 8).Value = 43000
$ws.Cells.Item(65, 10).Value = 43000
$ws.Cells.Item(65, 12).Value = 129000
$ws.Cells.Item(65, 14).Value = -135864
$ws.Cells.Item(76, 8).Value = 154380
$ws.Cells.Item(76, 10).Value = 154380
$ws.Cells.Item(76, 12).Value = 154380
$ws.Cells.Item(76, 14).Value = -155010
$ws.Cells.Item(79, 8).Value = 154380
$ws.Cells.Item(79, 10).Value = 154380
$ws.Cells.Item(79, 12).Value = 154380
$ws.Cells.Item(79, 14).Value = -156564
$ws.Cells.Item(86, 8).Value = 72942.5
$ws.Cells.Item(86, 9).Value = 59037.363
$ws.Cells.Item(86, 10).Value = 92062.06
$ws.Cells.Item(86, 11).Value = 59037.363
$ws.Cells.Item(86, 12).Value = 92062.06
$ws.Cells.Item(86, 13).Value = -57914.363
$ws.Cells.Item(86, 14).Value = -94308.06
$ws.Cells.Item(89, 8).Value = 72942.5
$ws.Cells.Item(89, 9).Value = 59037.363
$ws.Cells.Item(89, 10).Value = 92062.06
$ws.Cells.Item(89, 11).Value = 295186.815
$ws.Cells.Item(89, 12).Value = 460310.3
$ws.Cells.Item(89, 13).Value = -289570.815
$ws.Cells.Item(89, 14).Value = -471542.3
$ws.Cells.Item(99, 8).Value = 4788146.5
$ws.Cells.Item(99, 9).Value = 2522.7693
$ws.Cells.Item(99, 10).Value = 15156998
$ws.Cells.Item(99, 11).Value = 2522.7693
$ws.Cells.Item(99, 12).Value = 15156998
$ws.Cells.Item(99, 13).Value = -1024.7693
$ws.Cells.Item(99, 14).Value = -15159994
$ws.Cells.Item(107, 8).Value = 56253010
$ws.Cells.Item(107, 9).Value = 93752216
$ws.Cells.Item(107, 10).Value = 4193.125
$ws.Cells.Item(107, 11).Value = 93752216
$ws.Cells.Item(107, 12).Value = 4193.125
$ws.Cells.Item(107, 13).Value = -93750296
$ws.Cells.Item(107, 14).Value = -8033.125
$ws.Cells.Item(134, 8).Value = 4756.795
$ws.Cells.Item(134, 9).Value = 1827.9642
$ws.Cells.Item(134, 10).Value = 12212
$ws.Cells.Item(134, 11).Value = 5483.892599999999
$ws.Cells.Item(134, 12).Value = 36636
$ws.Cells.Item(134, 13).Value = -2948.892599999999
$ws.Cells.Item(134, 14).Value = -41706

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 3380.7778
$ws.Cells.Item(6, 10).Value = 5000
$ws.Cells.Item(6, 12).Value = 5000
$ws.Cells.Item(6, 14).Value = -5226
$ws.Cells.Item(10, 8).Value = 249.16667
$ws.Cells.Item(10, 9).Value = 249.16667
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 249.16667
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = -110.16667
$ws.Cells.Item(10, 14).Value = $null
$ws.Cells.Item(16, 8).Value = 5283.2104
$ws.Cells.Item(16, 9).Value = 4952.6665
$ws.Cells.Item(16, 10).Value = 5849.857
$ws.Cells.Item(16, 11).Value = 4952.6665
$ws.Cells.Item(16, 12).Value = 5849.857
$ws.Cells.Item(16, 13).Value = -4665.6665
$ws.Cells.Item(16, 14).Value = -6423.857
$ws.Cells.Item(18, 8).Value = 38561
$ws.Cells.Item(18, 10).Value = 38561
$ws.Cells.Item(18, 12).Value = 38561
$ws.Cells.Item(18, 14).Value = -39021
$ws.Cells.Item(28, 8).Value = 48397.25
$ws.Cells.Item(28, 10).Value = 48397.25
$ws.Cells.Item(28, 12).Value = 48397.25
$ws.Cells.Item(28, 14).Value = -48887.25
$ws.Cells.Item(31, 8).Value = 10072.521
$ws.Cells.Item(31, 9).Value = 4587.5454
$ws.Cells.Item(31, 10).Value = 15100.417
$ws.Cells.Item(31, 11).Value = 4587.5454
$ws.Cells.Item(31, 12).Value = 15100.417
$ws.Cells.Item(31, 13).Value = -4292.5454
$ws.Cells.Item(31, 14).Value = -15690.417
$ws.Cells.Item(34, 8).Value = 10072.521
$ws.Cells.Item(34, 9).Value = 4587.5454
$ws.Cells.Item(34, 10).Value = 15100.417
$ws.Cells.Item(34, 11).Value = 4587.5454
$ws.Cells.Item(34, 12).Value = 15100.417
$ws.Cells.Item(34, 13).Value = -4385.5454
$ws.Cells.Item(34, 14).Value = -15504.417
$ws.Cells.Item(35, 8).Value = 7550.1665
$ws.Cells.Item(35, 9).Value = 5506.727
$ws.Cells.Item(35, 10).Value = 30028
$ws.Cells.Item(35, 11).Value = 5506.727
$ws.Cells.Item(35, 12).Value = 30028
$ws.Cells.Item(35, 13).Value = -5212.727
$ws.Cells.Item(35, 14).Value = -30616
$ws.Cells.Item(38, 8).Value = 26071
$ws.Cells.Item(38, 10).Value = 28577.6
$ws.Cells.Item(38, 12).Value = 28577.6
$ws.Cells.Item(38, 14).Value = -29331.6
$ws.Cells.Item(46, 8).Value = 26071
$ws.Cells.Item(46, 10).Value = 28577.6
$ws.Cells.Item(46, 12).Value = 28577.6
$ws.Cells.Item(46, 14).Value = -28999.6
$ws.Cells.Item(53, 8).Value = 350266.75
$ws.Cells.Item(53, 10).Value = 350266.75
$ws.Cells.Item(53, 12).Value = 350266.75
$ws.Cells.Item(53, 14).Value = -351480.75
$ws.Cells.Item(94, 8).Value = 893.7727
$ws.Cells.Item(94, 10).Value = 754.82355
$ws.Cells.Item(94, 12).Value = 754.82355
$ws.Cells.Item(94, 14).Value = -1656.82355
$ws.Cells.Item(95, 8).Value = 176027.72
$ws.Cells.Item(95, 10).Value = 176027.72
$ws.Cells.Item(95, 12).Value = 176027.72
$ws.Cells.Item(95, 14).Value = -181519.72
$ws.Cells.Item(99, 8).Value = 7661.8335
$ws.Cells.Item(99, 9).Value = 1998.5
$ws.Cells.Item(99, 11).Value = 1998.5
$ws.Cells.Item(99, 13).Value = -500.5
$ws.Cells.Item(105, 9).Value = 7938707.5
$ws.Cells.Item(105, 11).Value = 7938707.5
$ws.Cells.Item(105, 13).Value = -7936960.5
$ws.Cells.Item(107, 8).Value = 989.1539
$ws.Cells.Item(107, 9).Value = 528.44446
$ws.Cells.Item(107, 11).Value = 528.44446
$ws.Cells.Item(107, 13).Value = 1391.55554
$ws.Cells.Item(108, 8).Value = 34844.75
$ws.Cells.Item(108, 9).Value = 19999
$ws.Cells.Item(108, 10).Value = 39793.332
$ws.Cells.Item(108, 11).Value = 19999
$ws.Cells.Item(108, 12).Value = 39793.332
$ws.Cells.Item(108, 13).Value = -16159
$ws.Cells.Item(108, 14).Value = -47473.332
$ws.Cells.Item(113, 8).Value = 5283.2104
$ws.Cells.Item(113, 9).Value = 4952.6665
$ws.Cells.Item(113, 10).Value = 5849.857
$ws.Cells.Item(113, 11).Value = 4952.6665
$ws.Cells.Item(113, 12).Value = 5849.857
$ws.Cells.Item(113, 13).Value = -2782.6665
$ws.Cells.Item(113, 14).Value = -10189.857
$ws.Cells.Item(122, 8).Value = 2087.6155
$ws.Cells.Item(122, 9).Value = 2212.5
$ws.Cells.Item(122, 11).Value = 6637.5
$ws.Cells.Item(122, 13).Value = -4187.5
$ws.Cells.Item(126, 8).Value = 7661.8335
$ws.Cells.Item(126, 9).Value = 1998.5
$ws.Cells.Item(126, 11).Value = 5995.5
$ws.Cells.Item(126, 13).Value = -3525.5
$ws.Cells.Item(132, 8).Value = 9003.091
$ws.Cells.Item(132, 9).Value = 2806.8
$ws.Cells.Item(132, 10).Value = 14166.667
$ws.Cells.Item(132, 11).Value = 8420.400000000001
$ws.Cells.Item(132, 12).Value = 42500.001
$ws.Cells.Item(132, 13).Value = -5890.400000000001
$ws.Cells.Item(132, 14).Value = -47560.001
$ws.Cells.Item(134, 8).Value = 3680.6956
$ws.Cells.Item(134, 9).Value = 1634.2
$ws.Cells.Item(134, 10).Value = 10192.272
$ws.Cells.Item(134, 11).Value = 4902.6
$ws.Cells.Item(134, 12).Value = 30576.816
$ws.Cells.Item(134, 13).Value = -2367.6
$ws.Cells.Item(134, 14).Value = -35646.81600000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2146.5
$ws.Cells.Item(5, 9).Value = 1913.1364
$ws.Cells.Item(5, 10).Value = 2467.375
$ws.Cells.Item(5, 11).Value = 5739.4092
$ws.Cells.Item(5, 12).Value = 7402.125
$ws.Cells.Item(5, 13).Value = -5627.4092
$ws.Cells.Item(5, 14).Value = -7626.125
$ws.Cells.Item(14, 8).Value = 11906622
$ws.Cells.Item(14, 9).Value = 11906622
$ws.Cells.Item(14, 11).Value = 35719866
$ws.Cells.Item(14, 13).Value = -35719693
$ws.Cells.Item(23, 8).Value = 71428760
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 71428760
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 214286280
$ws.Cells.Item(23, 13).Value = $null
$ws.Cells.Item(23, 14).Value = -214286750
$ws.Cells.Item(31, 8).Value = 433.33334
$ws.Cells.Item(31, 9).Value = 433.33334
$ws.Cells.Item(31, 11).Value = 1300.00002
$ws.Cells.Item(31, 13).Value = -1012.00002
$ws.Cells.Item(38, 8).Value = 627.8
$ws.Cells.Item(38, 9).Value = 585
$ws.Cells.Item(38, 10).Value = 692
$ws.Cells.Item(38, 11).Value = 1755
$ws.Cells.Item(38, 12).Value = 2076
$ws.Cells.Item(38, 13).Value = -1408
$ws.Cells.Item(38, 14).Value = -2770
$ws.Cells.Item(41, 8).Value = 1236.4117
$ws.Cells.Item(41, 10).Value = 1313
$ws.Cells.Item(41, 12).Value = 3939
$ws.Cells.Item(41, 14).Value = -4615
$ws.Cells.Item(68, 8).Value = 2217.3928
$ws.Cells.Item(68, 9).Value = 1832.3334
$ws.Cells.Item(68, 10).Value = 2322.4092
$ws.Cells.Item(68, 11).Value = 5497.0002
$ws.Cells.Item(68, 12).Value = 6967.2276
$ws.Cells.Item(68, 13).Value = -4686.0002
$ws.Cells.Item(68, 14).Value = -8589.2276
$ws.Cells.Item(71, 8).Value = 2217.3928
$ws.Cells.Item(71, 9).Value = 1832.3334
$ws.Cells.Item(71, 10).Value = 2322.4092
$ws.Cells.Item(71, 11).Value = 16491.0006
$ws.Cells.Item(71, 12).Value = 20901.6828
$ws.Cells.Item(71, 13).Value = -12435.0006
$ws.Cells.Item(71, 14).Value = -29013.6828
$ws.Cells.Item(80, 8).Value = 28575894
$ws.Cells.Item(80, 9).Value = 22731272
$ws.Cells.Item(80, 10).Value = 38466790
$ws.Cells.Item(80, 11).Value = 68193816
$ws.Cells.Item(80, 12).Value = 115400370
$ws.Cells.Item(80, 13).Value = -68192880
$ws.Cells.Item(80, 14).Value = -115402242
$ws.Cells.Item(83, 8).Value = 28575894
$ws.Cells.Item(83, 9).Value = 22731272
$ws.Cells.Item(83, 10).Value = 38466790
$ws.Cells.Item(83, 11).Value = 204581448
$ws.Cells.Item(83, 12).Value = 346201110
$ws.Cells.Item(83, 13).Value = -204576768
$ws.Cells.Item(83, 14).Value = -346210470
$ws.Cells.Item(92, 8).Value = 1145.2142
$ws.Cells.Item(92, 9).Value = 812.6667
$ws.Cells.Item(92, 10).Value = 1302.7368
$ws.Cells.Item(92, 11).Value = 2438.0001
$ws.Cells.Item(92, 12).Value = 3908.2104
$ws.Cells.Item(92, 13).Value = -1190.0001
$ws.Cells.Item(92, 14).Value = -6404.2104
$ws.Cells.Item(98, 8).Value = 3999.5
$ws.Cells.Item(98, 10).Value = 4999
$ws.Cells.Item(98, 12).Value = 14997
$ws.Cells.Item(98, 14).Value = -17993
$ws.Cells.Item(109, 8).Value = 2792.9
$ws.Cells.Item(109, 9).Value = 488.16666
$ws.Cells.Item(109, 10).Value = 6250
$ws.Cells.Item(109, 11).Value = 1464.49998
$ws.Cells.Item(109, 12).Value = 18750
$ws.Cells.Item(109, 13).Value = -424.4999800000001
$ws.Cells.Item(109, 14).Value = -20830
$ws.Cells.Item(122, 8).Value = 833057.75
$ws.Cells.Item(122, 9).Value = 2176827.8
$ws.Cells.Item(122, 10).Value = 1200.0952
$ws.Cells.Item(122, 11).Value = 19591450.2
$ws.Cells.Item(122, 12).Value = 10800.8568
$ws.Cells.Item(122, 13).Value = -19589000.2
$ws.Cells.Item(122, 14).Value = -15700.8568
$ws.Cells.Item(132, 8).Value = 5012.227
$ws.Cells.Item(132, 9).Value = 2704.8125
$ws.Cells.Item(132, 10).Value = 11165.333
$ws.Cells.Item(132, 11).Value = 24343.3125
$ws.Cells.Item(132, 12).Value = 100487.997
$ws.Cells.Item(132, 13).Value = -21813.3125
$ws.Cells.Item(132, 14).Value = -105547.997
$ws.Cells.Item(135, 8).Value = 2146.5
$ws.Cells.Item(135, 9).Value = 1913.1364
$ws.Cells.Item(135, 10).Value = 2467.375
$ws.Cells.Item(135, 11).Value = 17218.2276
$ws.Cells.Item(135, 12).Value = 22206.375
$ws.Cells.Item(135, 13).Value = -14683.2276
$ws.Cells.Item(135, 14).Value = -27276.375
$ws.Cells.Item(136, 8).Value = 3450.3333
$ws.Cells.Item(136, 9).Value = 3450.3333
$ws.Cells.Item(136, 11).Value = 10350.9999
$ws.Cells.Item(136, 13).Value = -5250.999899999999
$ws.Cells.Item(137, 8).Value = 85883.39999999999
$ws.Cells.Item(137, 10).Value = 113458.9
$ws.Cells.Item(137, 12).Value = 340376.7
$ws.Cells.Item(137, 14).Value = -350576.7
$ws.Cells.Item(139, 8).Value = 139559.73
$ws.Cells.Item(139, 9).Value = 168684.22
$ws.Cells.Item(139, 10).Value = 8499.5
$ws.Cells.Item(139, 11).Value = 506052.66
$ws.Cells.Item(139, 12).Value = 25498.5
$ws.Cells.Item(139, 13).Value = -500912.66
$ws.Cells.Item(139, 14).Value = -35778.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 3603.625
$ws.Cells.Item(43, 9).Value = 804.8333
$ws.Cells.Item(43, 11).Value = 804.8333
$ws.Cells.Item(43, 13).Value = -653.8333
$ws.Cells.Item(57, 8).Value = 51805.5
$ws.Cells.Item(57, 10).Value = 51805.5
$ws.Cells.Item(57, 12).Value = 51805.5
$ws.Cells.Item(57, 14).Value = -53445.5
$ws.Cells.Item(80, 8).Value = 2506.647
$ws.Cells.Item(80, 9).Value = 1875.2
$ws.Cells.Item(80, 10).Value = 3408.7144
$ws.Cells.Item(80, 11).Value = 1875.2
$ws.Cells.Item(80, 12).Value = 3408.7144
$ws.Cells.Item(80, 13).Value = -877.2
$ws.Cells.Item(80, 14).Value = -5404.7144
$ws.Cells.Item(83, 8).Value = 2506.647
$ws.Cells.Item(83, 9).Value = 1875.2
$ws.Cells.Item(83, 10).Value = 3408.7144
$ws.Cells.Item(83, 11).Value = 9376
$ws.Cells.Item(83, 12).Value = 17043.572
$ws.Cells.Item(83, 13).Value = -4384
$ws.Cells.Item(83, 14).Value = -27027.572
$ws.Cells.Item(97, 8).Value = 1667.5555
$ws.Cells.Item(97, 9).Value = 1739.8125
$ws.Cells.Item(97, 11).Value = 1739.8125
$ws.Cells.Item(97, 13).Value = -1243.8125
$ws.Cells.Item(102, 8).Value = 6732.7646
$ws.Cells.Item(102, 9).Value = 6389.7856
$ws.Cells.Item(102, 10).Value = 8333.333000000001
$ws.Cells.Item(102, 11).Value = 6389.7856
$ws.Cells.Item(102, 12).Value = 8333.333000000001
$ws.Cells.Item(102, 13).Value = -4767.7856
$ws.Cells.Item(102, 14).Value = -11577.333
$ws.Cells.Item(113, 8).Value = 4244.5557
$ws.Cells.Item(113, 9).Value = 2034.1
$ws.Cells.Item(113, 10).Value = 5544.8237
$ws.Cells.Item(113, 11).Value = 2034.1
$ws.Cells.Item(113, 12).Value = 5544.8237
$ws.Cells.Item(113, 13).Value = 135.9000000000001
$ws.Cells.Item(113, 14).Value = -9884.823700000001
$ws.Cells.Item(122, 8).Value = 33708.848
$ws.Cells.Item(122, 9).Value = 61559.293
$ws.Cells.Item(122, 10).Value = 4117.75
$ws.Cells.Item(122, 11).Value = 184677.879
$ws.Cells.Item(122, 12).Value = 12353.25
$ws.Cells.Item(122, 13).Value = -182227.879
$ws.Cells.Item(122, 14).Value = -17253.25
$ws.Cells.Item(124, 8).Value = 51274.5
$ws.Cells.Item(124, 10).Value = 51274.5
$ws.Cells.Item(124, 12).Value = 51274.5
$ws.Cells.Item(124, 14).Value = -61094.5
$ws.Cells.Item(126, 8).Value = 7368.5
$ws.Cells.Item(126, 9).Value = 6006
$ws.Cells.Item(126, 10).Value = 8049.75
$ws.Cells.Item(126, 11).Value = 18018
$ws.Cells.Item(126, 12).Value = 24149.25
$ws.Cells.Item(126, 13).Value = -15548
$ws.Cells.Item(126, 14).Value = -29089.25
$ws.Cells.Item(132, 8).Value = 4697.5
$ws.Cells.Item(132, 9).Value = 3058.0833
$ws.Cells.Item(132, 10).Value = 6336.9165
$ws.Cells.Item(132, 11).Value = 9174.249899999999
$ws.Cells.Item(132, 12).Value = 19010.7495
$ws.Cells.Item(132, 13).Value = -6644.249899999999
$ws.Cells.Item(132, 14).Value = -24070.7495
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).Value = $null
$ws.Cells.Item(136, 8).Value = 75455.875
$ws.Cells.Item(136, 10).Value = 75455.875
$ws.Cells.Item(136, 12).Value = 226367.625
$ws.Cells.Item(136, 14).Value = -231467.625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6070.737
$ws.Cells.Item(7, 9).Value = 3904.5293
$ws.Cells.Item(7, 10).Value = 7824.3335
$ws.Cells.Item(7, 11).Value = 3904.5293
$ws.Cells.Item(7, 12).Value = 7824.3335
$ws.Cells.Item(7, 13).Value = -3792.5293
$ws.Cells.Item(7, 14).Value = -8048.3335
$ws.Cells.Item(16, 8).Value = 3457
$ws.Cells.Item(16, 9).Value = 3389.125
$ws.Cells.Item(16, 11).Value = 3389.125
$ws.Cells.Item(16, 13).Value = -3219.125
$ws.Cells.Item(32, 8).Value = 13000
$ws.Cells.Item(32, 9).Value = 13000
$ws.Cells.Item(32, 11).Value = 13000
$ws.Cells.Item(32, 13).Value = -12683
$ws.Cells.Item(40, 8).Value = 45461244
$ws.Cells.Item(40, 9).Value = 71433670
$ws.Cells.Item(40, 11).Value = 71433670
$ws.Cells.Item(40, 13).Value = -71433534
$ws.Cells.Item(46, 8).Value = 3062.889
$ws.Cells.Item(46, 9).Value = 1521.3334
$ws.Cells.Item(46, 10).Value = 3833.6667
$ws.Cells.Item(46, 11).Value = 1521.3334
$ws.Cells.Item(46, 12).Value = 3833.6667
$ws.Cells.Item(46, 13).Value = -1333.3334
$ws.Cells.Item(46, 14).Value = -4209.6667
$ws.Cells.Item(68, 8).Value = 250025950
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).Value = $null
$ws.Cells.Item(69, 8).Value = 43000
$ws.Cells.Item(69, 10).Value = 43000
$ws.Cells.Item(69, 12).Value = 43000
$ws.Cells.Item(69, 14).Value = -44622
$ws.Cells.Item(71, 8).Value = 250025950
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).Value = $null
$ws.Cells.Item(72, 8).Value = 43000
$ws.Cells.Item(72, 10).Value = 43000
$ws.Cells.Item(72, 12).Value = 129000
$ws.Cells.Item(72, 14).Value = -137112
$ws.Cells.Item(87, 8).Value = 46595
$ws.Cells.Item(87, 10).Value = 46595
$ws.Cells.Item(87, 12).Value = 46595
$ws.Cells.Item(87, 14).Value = -48841
$ws.Cells.Item(90, 8).Value = 46595
$ws.Cells.Item(90, 10).Value = 46595
$ws.Cells.Item(90, 12).Value = 139785
$ws.Cells.Item(90, 14).Value = -151017
$ws.Cells.Item(95, 8).Value = 60000
$ws.Cells.Item(95, 10).Value = 60000
$ws.Cells.Item(95, 12).Value = 60000
$ws.Cells.Item(95, 14).Value = -65492
$ws.Cells.Item(100, 8).Value = 5666.3335
$ws.Cells.Item(100, 9).Value = 4997.5
$ws.Cells.Item(100, 11).Value = 4997.5
$ws.Cells.Item(100, 13).Value = -4456.5
$ws.Cells.Item(122, 8).Value = 5128.8223
$ws.Cells.Item(122, 9).Value = 4660.2
$ws.Cells.Item(122, 10).Value = 6066.067
$ws.Cells.Item(122, 11).Value = 13980.6
$ws.Cells.Item(122, 12).Value = 18198.201
$ws.Cells.Item(122, 13).Value = -11530.6
$ws.Cells.Item(122, 14).Value = -23098.201
$ws.Cells.Item(126, 8).Value = 6070.737
$ws.Cells.Item(126, 9).Value = 3904.5293
$ws.Cells.Item(126, 10).Value = 7824.3335
$ws.Cells.Item(126, 11).Value = 11713.5879
$ws.Cells.Item(126, 12).Value = 23473.0005
$ws.Cells.Item(126, 13).Value = -9243.5879
$ws.Cells.Item(126, 14).Value = -28413.0005
$ws.Cells.Item(132, 8).Value = 5813.3
$ws.Cells.Item(132, 9).Value = 3111.95
$ws.Cells.Item(132, 10).Value = 8514.65
$ws.Cells.Item(132, 11).Value = 9335.849999999999
$ws.Cells.Item(132, 12).Value = 25543.95
$ws.Cells.Item(132, 13).Value = -6805.849999999999
$ws.Cells.Item(132, 14).Value = -30603.95
$ws.Cells.Item(136, 8).Value = 10174.289
$ws.Cells.Item(136, 9).Value = 7444.0586
$ws.Cells.Item(136, 11).Value = 22332.1758
$ws.Cells.Item(136, 13).Value = -19782.1758

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(7, 8).Value = 1222
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 14).Value = $null
$ws.Cells.Item(9, 8).Value = 60000000
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 13).Value = $null
$ws.Cells.Item(81, 8).Value = 17506516
$ws.Cells.Item(81, 9).Value = 1112687.5
$ws.Cells.Item(81, 10).Value = 66688000
$ws.Cells.Item(81, 11).Value = 2225375
$ws.Cells.Item(81, 12).Value = 133376000
$ws.Cells.Item(81, 13).Value = -2224314
$ws.Cells.Item(81, 14).Value = -133378122
$ws.Cells.Item(84, 8).Value = 17506516
$ws.Cells.Item(84, 9).Value = 1112687.5
$ws.Cells.Item(84, 10).Value = 66688000
$ws.Cells.Item(84, 11).Value = 11126875
$ws.Cells.Item(84, 12).Value = 666880000
$ws.Cells.Item(84, 13).Value = -11121571
$ws.Cells.Item(84, 14).Value = -666890608
$ws.Cells.Item(95, 8).Value = 53952
$ws.Cells.Item(95, 10).Value = 53952
$ws.Cells.Item(95, 12).Value = 53952
$ws.Cells.Item(95, 14).Value = -59444
$ws.Cells.Item(96, 8).Value = 3158
$ws.Cells.Item(96, 9).Value = 2688.4443
$ws.Cells.Item(96, 11).Value = 2688.4443
$ws.Cells.Item(96, 13).Value = -1315.4443
$ws.Cells.Item(97, 8).Value = 33572
$ws.Cells.Item(97, 10).Value = 33572
$ws.Cells.Item(97, 12).Value = 33572
$ws.Cells.Item(97, 14).Value = -35554
$ws.Cells.Item(98, 8).Value = 333368260
$ws.Cells.Item(98, 10).Value = 333368260
$ws.Cells.Item(98, 12).Value = 333368260
$ws.Cells.Item(98, 14).Value = -333374250
$ws.Cells.Item(105, 8).Value = 67402
$ws.Cells.Item(105, 10).Value = 67402
$ws.Cells.Item(105, 12).Value = 67402
$ws.Cells.Item(105, 14).Value = -74390
$ws.Cells.Item(113, 8).Value = 11437.565
$ws.Cells.Item(113, 9).Value = 18203.428
$ws.Cells.Item(113, 11).Value = 54610.284
$ws.Cells.Item(113, 13).Value = -52440.284
$ws.Cells.Item(122, 8).Value = 7205177.5
$ws.Cells.Item(122, 9).Value = 10289076
$ws.Cells.Item(122, 10).Value = 9414.143
$ws.Cells.Item(122, 11).Value = 30867228
$ws.Cells.Item(122, 12).Value = 28242.429
$ws.Cells.Item(122, 13).Value = -30864778
$ws.Cells.Item(122, 14).Value = -33142.429
$ws.Cells.Item(125, 8).Value = 50830.4
$ws.Cells.Item(125, 10).Value = 50830.4
$ws.Cells.Item(125, 12).Value = 50830.4
$ws.Cells.Item(125, 14).Value = -60670.4
$ws.Cells.Item(126, 8).Value = 3636.52
$ws.Cells.Item(126, 9).Value = 1382
$ws.Cells.Item(126, 11).Value = 4146
$ws.Cells.Item(126, 13).Value = -1676
$ws.Cells.Item(132, 8).Value = 12043.326
$ws.Cells.Item(132, 9).Value = 7872.273
$ws.Cells.Item(132, 10).Value = 20646.125
$ws.Cells.Item(132, 11).Value = 23616.819
$ws.Cells.Item(132, 12).Value = 61938.375
$ws.Cells.Item(132, 13).Value = -21086.819
$ws.Cells.Item(132, 14).Value = -66998.375
$ws.Cells.Item(136, 8).Value = 35376.5
$ws.Cells.Item(136, 9).Value = 1667.3
$ws.Cells.Item(136, 10).Value = 91558.5
$ws.Cells.Item(136, 11).Value = 5001.9
$ws.Cells.Item(136, 12).Value = 274675.5
$ws.Cells.Item(136, 13).Value = -2451.9
$ws.Cells.Item(136, 14).Value = -279775.5

Write-Host "Applied 708 cell updates and 8 cell clears"
